$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rework the "charge item card" / monster-flip / trap-move / loot style
#     event cards to cost "1 time" instead of only attribute/hand discard,
#     and rename leveling up into "skill points" ---------------------------

$ws.Range("D2").Value  = "消耗1时间，或使用1任意属性，或将1张手牌洗回主牌堆：重抽房间区所有牌。使用1张《感知》发动本牌时，可以再重抽任意张手牌。"
$ws.Range("D3").Value  = "消耗1时间，或使用1任意属性，或将1张手牌洗回主牌堆：将房间区最前方的所有怪物牌横置。使用1张《敏捷》发动本牌时，可以将房间区所有怪物牌横置。"
$ws.Range("D4").Value  = "消耗1时间，或使用1任意属性，或将1张手牌洗回主牌堆：将房间区最前方任意1张牌移动到战场敌人列第一行。使用1张《敏捷》发动本牌时，可以将房间区任意1张牌移动到战场敌人列第一行。"
$ws.Range("D5").Value  = "消耗1时间，或使用1任意属性，或将1张手牌洗回主牌堆：获得遭遇牌堆第1张战利品牌。使用1张《感知》发动本牌时，转而翻开遭遇牌堆前3张战利品牌中选1张获得。"
$ws.Range("D6").Value  = "消耗1时间，或使用1任意属性，或将1张手牌洗回主牌堆：选房间区1张陷阱牌移动到房间区任意位置。使用1张《感知》发动本牌时，可以转而选弃牌堆1张陷阱牌放到房间区任意位置。"
$ws.Range("D7").Value  = "消耗1时间，或使用1任意属性，或将1张手牌洗回主牌堆：从遭遇牌堆翻开3张牌，获得其中的战利品牌。使用1张《敏捷》发动本牌时，可以额外翻开2张牌。"
$ws.Range("D8").Value  = "消耗2时间：将主牌堆第1张怪物牌放在房间区任意非空列顶端，然后获得遭遇牌堆第1张战利品牌，再获得遗物牌堆顶的1张遗物牌。"
$ws.Range("D9").Value  = "多选，可重复：①支付1金币，获得1道具点。②支付2道具点，获得1金币。②支付3金币，从遗物牌堆翻开3张牌，选其中1张获得。"
$ws.Range("D10").Value = "消耗5时间，然后使用1任意属性或将1张手牌洗回主牌堆：获得1技能点。使用1张《智力》发动本牌时，可以少消耗2时间。"

# --- Add a new "祭坛" (Altar) card as row 11, introducing the buy cost for
#     upgrade cards (pay a resource to gain 1 item point) ------------------

$ws.Range("A11").Value = "祭坛"
$ws.Range("C11").Value = 3
$ws.Range("D11").Value = "将1张手牌送墓，或弃置1张战利品牌，或受到1点伤害：获得1道具点。"

# --- Update the view selection to match the saved session -----------------

$ws.Range("E9").Select()
